$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the affected cells to remain plain text so Excel does not
# auto-convert numeric-looking or date-looking strings (applied one
# address at a time - a combined multi-area range string does not
# reliably propagate the NumberFormat to every cell).
$textCellAddrs = @("I2","J2","K2","M2","O2","P2","R2","T2","U2","W2","X2","AA2","AC2","AD2")
foreach ($addr in $textCellAddrs) {
    $ws.Range($addr).NumberFormat = "@"
}

# Update row 2 cells with the new trade values
$ws.Range("I2").Value = "2022-04-21"
$ws.Range("J2").Value = "NIFTY2242116850PE"
$ws.Range("K2").Value = "33"
$ws.Range("M2").Value = "29"
$ws.Range("O2").Value = "Percentage"
$ws.Range("P2").Value = "7"
$ws.Range("R2").Value = "2"
$ws.Range("T2").Value = "15"
$ws.Range("U2").Value = "3"
$ws.Range("W2").Value = "existing"
$ws.Range("X2").Value = "NO"
$ws.Range("AA2").Value = "new"
$ws.Range("AC2").Value = "YES"
$ws.Range("AD2").Value = "existing"

# Remove the now-duplicate row 3 entirely (shrinks dimension from AE3 to AE2)
$ws.Rows.Item(3).Delete()
